# PrecioFrutaHortalizas / Coliflor weekly update:
# A new weekly record is inserted at row 1102 (D=45194, Primera, 2000/700/800/750/750),
# which pushes all subsequent records (rows 1102-1218) down by one row to 1103-1219.
# The used range grows from A1:R1218 to A1:R1219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the block of existing records (A1102:R1218) down one row, into A1103:R1219.
# This carries every column (dates, quality, volumes, prices, etc.) for each record
# down to the next row, making room for the new record at row 1102.
$srcRange = $ws.Range("A1102:R1218")
$srcValues = $srcRange.Value2
$dstRange = $ws.Range("A1103:R1219")
$dstRange.Value2 = $srcValues

# Carry over number formatting / styles (e.g. the date style on column D) the same way,
# since the new row 1219 doesn't have any formatting of its own yet.
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

# Now overwrite row 1102 with the new weekly record.
$ws.Range("A1102").Value2 = 8
$ws.Range("B1102").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C1102").Value2 = "Coquimbo"
$ws.Range("D1102").Value2 = 45194
$ws.Range("E1102").Value2 = 4
$ws.Range("F1102").Value2 = 100112008
$ws.Range("G1102").Value2 = "Coliflor"
$ws.Range("H1102").Value2 = "Sin especificar"
$ws.Range("I1102").Value2 = "Primera"
$ws.Range("J1102").Value2 = 2000
$ws.Range("K1102").Value2 = 700
$ws.Range("L1102").Value2 = 800
$ws.Range("M1102").Value2 = 750
$ws.Range("N1102").Value2 = "`$/unidad"
$ws.Range("O1102").Value2 = "Provincia del Elquí"
$ws.Range("P1102").Value2 = 750
$ws.Range("Q1102").Value2 = 1
$ws.Range("R1102").Value2 = "Hortaliza"

Write-Host "done"
